$d = $word.ActiveDocument

# 1. Merge "spacetime" run back into surrounding text (removes spell-check proofErr wrapping)
$d.Content.Find.Execute("The Lattice-Field Medium (LFM) proposes that spacetime itself is a discrete, dynamical lattice of locally interacting cells.", $true, $false, $false, $false, $false, $true, 1, $false, "The Lattice-Field Medium (LFM) proposes that spacetime itself is a discrete, dynamical lattice of locally interacting cells.", 2)

# 2. behaviours -> behaviors (span across proofErr tags so they get merged away)
$d.Content.Find.Execute("emergent behaviours of", $true, $false, $false, $false, $false, $true, 1, $false, "emergent behaviors of", 2)

# 3. Merge "Variational" back into surrounding text
$d.Content.Find.Execute("5. **Variational gravity (Sec. 13)**: a dynamic", $true, $false, $false, $false, $false, $true, 1, $false, "5. **Variational gravity (Sec. 13)**: a dynamic", 2)

# 4. behaviour -> behavior (span across proofErr tags so they get merged away)
$d.Content.Find.Execute("relativistic behaviour are", $true, $false, $false, $false, $false, $true, 1, $false, "relativistic behavior are", 2)
